$wb = $excel.ActiveWorkbook

# --- Sheet 1: Circulars ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range('A2').Value = 'SEBI'
$ws1.Range('B2').Value = 'Circulars'
$ws1.Range('C2').Value = 2025
$ws1.Range('D2').Value = 'October'
$ws1.Range('E2').NumberFormat = '@'
$ws1.Range('E2').Value = '2025-10-13'
$ws1.Range('E2').ClearFormats()
$ws1.Range('F2').Value = 'Minimum information to be provided to the Audit Committee and Shareholders for approval of Related Party Transactions'
$ws1.Range('G2').Value = 'https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1760356560260.pdf'
$ws1.Range('H2').Value = '1760356560260.pdf'
$ws1.Range('I2').Value = '/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Circulars/2025/October/1760356560260.pdf'
$ws1.Range('J2').Value = 'NA'
$ws1.Range('K2').Value = '"SEBI has relaxed the requirements for approving related party transactions by listed entities from September 12, 2025. This change aims to facilitate ease of doing business and applies to transactions below 1% of annual consolidated turnover or INR 10 crores. Specific information requirements have been outlined in SEBI''s recent circular, available on their website. Listed entities must comply with the specified format and RPT Industry Standards for certain sections of LODR Regulations. The circular also outlines additional details required for shareholder approval of related party transactions."'
$ws1.Range('L2').Value = 'SEBI has approved a proposal for relaxation in the minimum information required for approval of related party transactions by listed entities, effective from September 12, 2025. This change aims to facilitate ease of doing business. The new requirements only apply to transactions that do not exceed 1% of annual consolidated turnover or INR 10 crores (whichever is lower), and specific information requirements have been outlined in SEBI''s Circular dated June 26, 2025 and Annexure-13A.
This section details an exemption for related party transactions not exceeding Rs. One Crore during a financial year, as per the RPT Industry Standards. From the effective date of this circular, listed entities must comply with the specified format and RPT Industry Standards for Part A and Part B of Section III-B of LODR Regulations. This circular takes immediate effect and is available on SEBI''s website. The document also outlines minimum information requirements for audit committee approval and shareholder approval of related party transactions.
Shareholders are requested to approve proposed Related Party Transactions (RPT). Alongside the requirements under the Companies Act 2013, the notice must provide additional details: management''s summary, justification for transaction benefits, specifics on loans/investments, valuation report availability, counter-party turnover representation, and other relevant information. No deadlines are specified in this text.'
# Remove the now-obsolete rows 3 and 4 (shrinks dimension to A1:L2)
$ws1.Rows('3:4').Delete()

# --- Sheet 3: was 'Consulatation Paper' -> 'Master circular' ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = 'Master circular'
$ws3.Range('A2').Value = 'SEBI'
$ws3.Range('B2').Value = 'Master circular'
$ws3.Range('C2').Value = 2025
$ws3.Range('D2').Value = 'October'
$ws3.Range('E2').NumberFormat = '@'
$ws3.Range('E2').Value = '2025-10-15'
$ws3.Range('E2').ClearFormats()
$ws3.Range('F2').Value = 'Master Circular for issue and listing of Non-convertible Securities, Securitised Debt Instruments, Security Receipts, Municipal Debt Securities and Commercial Paper'
$ws3.Range('G2').Value = 'https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1760532257519.pdf'
$ws3.Range('H2').Value = '1760532257519.pdf'
$ws3.Range('I2').Value = '/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Master circular/2025/October/1760532257519.pdf'
$ws3.Range('J2').Value = 'No.
CIR/MIRSD/4/2011, No
0, No
0, No.
1, No. 1, No
49, No. 1, No. 1, No.
9'
$ws3.Range('K2').Value = '"SEBI has issued a consolidated document updating regulations for various types of securities and debt instruments, affecting multiple financial entities such as issuers, stock exchanges, depositories, credit rating agencies, bankers, and intermediaries. Key changes include new regulations for green debt securities to combat greenwashing, provisions for the ASBA mechanism in public issues, and modifications for bond trading platforms. The document also introduces a Liquidity Window facility for investors in debt securities and updates on nominee directors, contributions to the Settlement Guarantee Fund, Legal Entity Identifier (LEI) for issuers, and Online Bond Platform Providers (OBPPs). Deadlines are provided in specific annexures."'
$ws3.Range('L2').Value = 'On October 15, 2025, a consolidated document was issued to update and streamline regulations for Non-convertible Securities, Securitised Debt Instruments, Security Receipts, Municipal Debt Securities, and Commercial Paper. The update incorporates provisions from multiple circulars issued till June 30, 2025, effectively rescinding earlier ones in this area. The affected parties include various financial entities such as issuers, stock exchanges, depositories, credit rating agencies, and bankers. All relevant entities are directed to disseminate the provisions on their websites, comply with the conditions laid out, and implement the necessary systems for adherence.
The regulatory document directs Recognized Stock Exchanges, Depositories, intermediaries, issuers, and other stakeholders to disseminate the provisions of the circular on their websites, comply with its conditions, implement necessary systems, make changes to bye-laws if needed, communicate with stakeholders, and monitor compliance. The directive is issued under various sections of the Act, 1992, SEBI regulations from 2008 to 2021, to protect investor interests and promote securities market development. This circular can be found on the SEBI website.
Chapter IX introduces new regulations for green debt securities to combat greenwashing, affecting issuers of such securities. Deadlines or critical dates are detailed in Annexure-IX-C-A (on page 94) and Annexure-IX-C-B (on page 97). Other chapters cover various aspects of debt securities, fundraising, reporting, and listing (Chapters X-XVIII), with additional details provided in numerous annexures.
This document introduces several regulatory updates and frameworks affecting Online Bond Platform Providers (OBPPs), trade execution platforms for certain securities, nominee directors, contributions to the Settlement Guarantee Fund, Legal Entity Identifier (LEI) for issuers, a Liquidity Window facility for investors in debt securities, and several abbreviations related to finance and securities. Deadlines or critical dates are not specified within this provided section.
Regulatory Update: The provisions outlined in Chapter I apply to public issues of securities under SEBI NCS Regulations (2021), ILDM Regulations (2015), and SDI Regulations (2008). Investors can utilize the ASBA mechanism by writing bank account numbers or mentioning UPI ID for payments up to Rs. 5 lakh per application. Applications can be submitted through SCSBs or intermediaries, with funds being blocked upon bid submission on the stock exchange platform.
Starting November 1, 2024 (per circular dated October 18, 2024), individual investors can submit bid-cum-application forms through various methods: intermediaries, SCSBs, or stock exchange platforms. For applications up to Rs. 5 Lakh, UPI is mandatory for funding blocking via intermediaries and optional for SCSBs and stock exchange platforms. All applications up to Rs. 5 Lakh through intermediaries must use UPI. Investors can also submit forms online using the linked online trading, demat, and bank account facility. SCSBs and intermediaries are responsible for uploading forms to the stock exchange bidding platform, providing acknowledgements, and ensuring sufficient, clearly demarcated funds for ASBA applications. SCSBs must have a separate account in another SEBI-registered SCSB for applying in public issues using the ASBA facility on their own account.
- SEBI advises SCSBs to have separate accounts for ASBA applications and ensure available funds. (3.7)
- Stock exchanges allow modification of bid details, provide real-time validation, update demand data, and facilitate status checking, alerts, and operational procedures. (4.1, 4.2, 7.1, 7.4)
- UPI, an instant payment system by NPCI, allows instant money transfer between bank accounts using a unique payment address. A Sponsor Bank is appointed to push mandate collect requests and/or payment instructions of investors into the UPI. (5.2, 5.3)
- The details of investor and bid entered on the stock exchange platform are validated with depositories in real time. Intermediaries provide guidance on UPI usage, and disclose commission, processing fees, and operational procedures. (6.1, 7.2, 7.3)
- Character length for each field of the schedule forwarded by intermediaries is uniformly prescribed by stock exchanges, following a specific format. (8)
- The bidding and validation process requires investors to have/create a UPI ID before application submission. (9.1a)
The document outlines updates for UPI (Unified Payments Interface) processes in relation to stock exchange bidding. Investors are now required to have a UPI ID before submitting applications, with bid details including linked UPI ID submitted via intermediaries or stock exchange apps/websites. Upon validation of PAN and demat account details, investors will receive an SMS about their application submission.
Post-validation, sponsor banks will request investors to authorize blocking of funds equivalent to the bid amount. This mandate must be accepted by 5 pm on the third working day following the bidding day, except for the last day when it should be done by 5 pm the next working day. Investors can modify certain details but not both DP ID or Client ID and PAN simultaneously. Mismatch bids will be processed till 1pm on T+1 day. The facility to reinitiate/resend UPI mandates is available only until 5 pm on the day of bidding.
- Until 5 pm on the day of bidding, investors can re-initiate or resend UPI mandates for block requests. (h)
- Upon successful validation of block requests, funds equivalent to the application amount will be blocked in investors'' accounts, and status updates will be shared with sponsor banks, stock exchanges, and intermediaries for reconciliation post issue closure. (i-k, 10.1-10.4)
- The allotment of securities occurs within five working days after issue closure, followed by credit to investors'' accounts and unblocking of excess funds. (10.2, 10.5)
- Issuers use an online platform provided by stock exchanges to receive applications for public debt securities and maintain a single escrow account. Registrars have an interface with the stock exchange platform for updates. (11.1a-c, 11.2a)
The section outlines key actions related to the issuance of non-convertible securities, municipal debt securities, and securitized debt instruments. Stock exchanges are required to provide an online platform for applications, including UPI mode, with a cap of INR 5 lakh. Intermediaries are responsible for accurate data entry and addressing investor grievances. The issuer shall appoint a sponsor bank as a conduit between stock exchanges and NPCI. All activities from issue closure to listing must be completed within T+6 working days, as detailed in Table 2.
Regulatory update: Listing timeline for public issues of debt securities and NCRPS reduced from T+6 to T+3 working days, effective immediately for one year, and permanently thereafter. This change aims to provide faster access to funds for issuers and investors. The revised timeline is optional during the transition period, with regulation 37 (2) of NCS Regulations applying after T+6 only in cases where T+3 was chosen but not met.
The provisions of regulation 37 (2) will only apply after T+6 working days during the voluntary use of the T+3 listing timeline, starting from November 01, 2024. This applies to public issues of debt securities and NCRPS. Mandatory application begins on November 01, 2025. The T+3 timeline for listing must be disclosed in Offer Documents for all public issues. An indicative timeline for activities related to listing is provided, detailing various stages from bid modification to trading commencement.
Regulatory update: The amendment requires investors, including residents and NRIs, to submit listing applications for public issues to stock exchanges. Stock exchanges will grant trading permission upon application. Trading commences on T+1 day, with specific deadlines for initiation, completion, and fund transfer by T+2 day and T+3 day. The application form details are provided in Annexure IIA and IIB for residents and NRIs respectively. Certain categories of persons, such as minors without a guardian, foreign entities, and those ineligible to contract under applicable statutory/regulatory requirements, are not eligible to participate in the issue. Application forms must be completed in full, in English, and errors in data entry are the responsibility of the applicant. Minimum application quantities are specified in the Prospectus.
Regulatory update: Changes to application procedures for Bond issues. Affects applicants, with requirements including use of specified Bonds in applications, attestation of non-Indian language thumb impressions and signatures by Magistrates or Notaries Public, use of valid and active Demat ID, Client ID, and PAN, submission of Applications in single or joint names (with specific orders), use of ASBA Accounts for electronic applications, and rejection of applications without PAN numbers. Critical dates not specified. Applicants must ensure eligibility as per terms of the Shelf & Tranche with Prospectus and applicable laws. Physical application submission to be made through SCSBs or Stock Exchange members. Application Forms must bear relevant stamps and contain correct ASBA Account numbers and UPI IDs. Revision/modification of applications during Issue Period is allowed, with original upload date considered for Allotment purposes.
Regulatory update for bond applications: HUFs applying through Karta must specify their name in the Application Form as ''XYZ Hindu Undivided Family applying through PQR''. All applicants should select the Series of Bonds they wish to apply for. ASBA Applicants need to provide correct details of their ASBA Account and submit forms at designated branches or with syndicate members, not directly to banks or the company.
Rejected applications include those without full payment, incorrect PAN, partnership firm submissions under individual names, applications by persons not competent to contract, GIR number instead of PAN, applications by OCBs, insufficient funds in ASBA Account, applications prohibited by SEBI or other authorities, untimely UPI Mandate approval, incomplete documents, incorrect signatures, unmatched Depository records, and applications uploaded after the Issue Closing.
Key regulatory action: The document outlines restrictions on applications for Bond issues, affecting various applicants such as those prohibited by SEBI, those with suspended demat accounts, and debarred persons from the capital markets. It also specifies methods for allotment based on the date of application, with preference given to first-come, first-served, and proportional allotment in case of oversubscription.
Critical dates: The Issue Closing Date is mentioned, with extensions possible by stock exchanges as applicable. The deadline for applicants to deliver their Application Forms is prescribed, along with instructions regarding the ASBA process. In case of oversubscription, the methodology for allotment includes a cutoff date one day prior to oversubscription.
Regulatory update: In case of equal entitlement among multiple applicants, a draw of lots will determine the basis of allotment (Section iii). If an applicant has applied for multiple Series of Bonds, allocation will be proportional to each applied Series (Section e(i)). The Company retains discretion to close the Issue early (Section iv). Investors can withdraw their applications at any time. Minimum subscription is 75%, and if not met, funds are unblocked with a potential interest penalty (Terms of the Issue, Point 1). The abridged prospectus format is specified in Annexure IID of the document, with certain formatting requirements (Chapter IIB, Section 1).
The revised format for Abridged Prospectus will apply to all public issues from October 1, 2023. Key changes include using larger font sizes, specific formatting for heads and sub-heads, tabular formats, and pointers for efficient understanding. Important information not included in designated sections can be added under ''any other information''. QR codes will be added to the documents for easy access to the prospectus. The format is detailed in Annexure – IID of the document. Cash flow conventions, disclosure of cash flows, and other related disclosures in the offer document must follow specific guidelines as per Regulation 14, Clause 3.3.34 Schedule I SEBI NCS Regulations, 2021.
1. For non-convertible securities, coupon/dividend payments falling on Sundays or holidays will be made on the next working day, while future payment dates remain as per original schedule. (2020-25)
2. In case of debt securities maturity, redemption proceeds will be paid on the previous working day if the maturity date falls on a Sunday or holiday. (No specific dates mentioned)
3. A uniform methodology for interest/dividend calculations in leap years is established, with 366 days considered as the denominator. (Applicable throughout tenor of security)
4. In case of a public issue, additional disclosures are required for provisions related to fictitious applications, board declarations, reservations, utilization details, promoter/director benefits, and material contracts. (Applicable for public issues)
5. For NBFCs, HFCs, or PFIs, specific disclosures are needed when the objects of a public issue involve loans to group companies. (Not specified)
Section IV of the SEBI NCS Regulations, 2021 introduces additional disclosure requirements for Non-Banking Finance Companies (NBFCs), Housing Finance Companies (HFCs), and Public Financial Institutions (PFI) in public issues. Key disclosures include details about loan policies, classifications of loans, exposures to top borrowers, and non-performing assets. NBFCs/HFCs are also required to provide portfolio summaries and disclose off-balance sheet items. These disclosures are to be made in a prescribed format and are based on draft templates for loan classifications, sectoral exposure, and loan ticket sizes. No specific deadlines were mentioned in the provided text.
Section 4.4 introduces a requirement for loan denomination based on ticket size at origination, affecting all loans and borrowers. No specific date or deadline is mentioned.
Section 4.5 involves classifying borrowers geographically, focusing on top 5 states with the highest AUM percentage. Again, no critical dates are provided.
In Section 4.6, details of loans overdue and classified as non-performing according to RBI''s stipulations are presented. Tables show movement of gross NPA (Non-Performing Assets) and provisions for NPA across various time periods. No specific dates or deadlines are mentioned.
Section 4.7 provides segment-wise gross NPA, differentiating between retail and wholesale sectors, with details on A Mortgages, Gold loans, Vehicle finance, MFI, MSME, Real estate, Promoter funding, and others. No critical dates or deadlines are provided.
Section 4.8 discusses the residual maturity profile of assets and liabilities, following RBI''s format and including categories such as deposit, advances, investment, borrowing, and FCA/FCL (Foreign Currency Assets/Liabilities). No specific dates or deadlines are mentioned.
Chapter V covers the denomination of issuance and trading of Non-convertible Securities, introducing changes for private placement of debt securities and non-convertible redeemable preference shares with a face value of either Rs. One lakh (for general cases) or Rs. Ten lakhs (for shelf placements from January 1, 2023). Various conditions and credit enhancements are also outlined in this section for securities issued on private placement basis at a face value of Rs. Ten Thousand. No critical dates or deadlines are mentioned.
This section outlines changes to credit enhancements allowed in certain securities, effective July 03, 2024, as per a circular dated the same day. Acceptable enhancements include guaranteed bonds, partially guaranteed bonds, SBLC-backed securities, debt backed by pledge of shares or other assets, guaranteed pooled bond issuance (PBI), obligor/co-obligor structures, cross default guarantee structures, and debt backed by payment waterfall/escrow with full guarantee from a third party.
For these enhanced securities, Credit Rating Agencies must verify that support is unconditional, irrevocable, and legally enforceable till all obligations are paid to investors, and the support provider has a lower probability of default compared to the issuer until the instruments are outstanding.
Additionally, for a shelf placement memorandum or General Information Document valid as of the circular''s effective date, issuers may now raise funds through tranche placement memorandum or Key Information Document at a face value of Rs. Ten Thousand, provided at least one Merchant Banker is appointed to carry out due diligence.
Effective January 1, 2023 (unless otherwise stated), eligible participants on the EBP Platform include Qualified Institutional Buyers (QIBs) as per SEBI ICDR Regulations, 2018 and authorized non-QIBs by issuers. Issues of securities made through the platform now include:
- Debt securities and NCRPS for private placements, shelf issues, or subsequent issues exceeding Rs. 20 crore in a financial year (previously Rs. 50 crore), per SEBI NCS Regulations, 2021.
- Issues of debt securities and NCRPS by issuers existing for less than three years, regardless of issue size, as per Clause 3.3.10 c. of Schedule I to the SEBI NCS Regulations, 2021.
- Specific non-equity instruments issued by banks/NBFCs listed under Chapter V of the SEBI NCS Regulations, 2021, regardless of issue size.
- Private placement of securitized debt instruments, security receipts, CPs, CDs, REIT units, SM REIT units, and InvIT units by issuers desiring to access the EBP platform.
- Issues of debt securities, NCRPS, and municipal debt securities on private placement basis below Rs. 20 crore may also utilize the EBP Platform.
Issuers are obligated to comply with relevant laws, provide Placement Memorandum and term sheet to the platform at least two working days prior to issue opening date (three days for first-time issuers), and disclose details of issue size, green shoe portion, etc., in these documents.
The regulatory document outlines changes effective May 16, 2025, affecting private placements of municipal debt securities or CPs/CDs. Issuers may now access the Electronic Book Building (EBP) platform for such transactions, regardless of issue size. Green shoe options are capped at five times the base issue size. Key details disclosed in the Placement Memorandum and term sheet include size, interest rate parameters, bid opening/closing dates, minimum bid lot, bidding method, allotment method, settlement method, and settlement cycle. Issuers must also disclose estimated cut-off yield at least one hour prior to bidding and provide post-issue details on the EBP platform. Participants are required to enroll with EBP for one-time access, with KYC verification handled by various parties depending on investor type. All eligible participants should have access to issue-specific information.
Key regulatory action: Changes in bidding process for investors on Electronic Book Building (EBP) platform.
Affects: Non-QIB investors, QIB investors, arrangers, and EBPs.
Critical dates or deadlines: Bids must be placed between 9 a.m. to 5 p.m., one day before bidding announcement, with a minimum of one hour bidding window. Issuer can provide details of eligible participants not later than one hour before the bidding start time. Changes in bidding date or time are allowed for a maximum of two times.
Eligible participants must confirm they do not use unfair access tools and EBPs must ensure no preferential access to any bidder. Bids cannot exceed Rs.100 crore or 5% of the base issue size (whichever is lower), with foreign portfolio investors able to bid through their custodians.
Arrangers can bid on behalf of multiple participants, but must disclose details about proprietary bids, client bids, and consolidated bids to EBP at the time of bidding. Bidding takes place on an anonymous order-driven system, with bids made in price, coupon (in %), spread in basis points (bps), or a combination thereof. Modification or cancellation of bids is allowed during the bidding period.
Regulatory update: The document outlines changes to bidding procedures for investors. Bidders can place multiple bids with coupon/ spread during the bidding period, which can be cancelled or modified only during the bidding period (except in the last 10 minutes). In the final 10 minutes, only downward revisions of coupon/ spread or upward modifications of price, or upward revisions in bid size are allowed. All bids have an audit trail and are disclosed on the EBP platform after bidding. The allotment and settlement amount for bidders are based on various factors including coupon specified by the issuer, coupon discovered during bidding, and pro-rata basis if multiple bids are at the cut-off coupon/ price/ spread. Anchor portions within the base issue size can be allocated by the issuer at their discretion, with limits set based on instrument rating (not exceeding 30% for AAA, 40% for AA+/AA, 50% for A+/A- and others). There is no bidding for anchor portion on the EBP platform.
Regulatory update: The maximum allocation to anchor investors in an issue is limited to 30% (as per May 16, 2025 circular). Anchor investors must provide electronic confirmation on EBP platform T-1 day prior to settlement. Details of anchor investors and their corresponding quantum allocated should be disclosed in the placement memorandum and term sheet. Remaining portion of the issue is open for bidding by eligible participants on EBP platform, with anchor investors allowed to participate if identified as such. Pay-in obligations must be fulfilled within T+1 or T+2 days, with non-compliance resulting in temporary access ban for 30 days (for allottees and anchor investors) or 7 days (for arrangers).
9.4: An arranger will be debarred from EBP platforms for seven days following three unpaid obligations on any Exchange-Based Platform (EBP).
9.5-9.6: Pay-ins must be made to an issuer''s escrow bank account or through clearing corporations of stock exchanges, with confirmation provided by the escrow bank and reconciliation performed by a Registrar to an Issue. Any pay-in from unapproved accounts will result in bid cancellation and investor debarment for 30 days.
10: An issuer may withdraw from any issue process at any time, resulting in a seven-day access restriction to EBP platforms. However, this restriction does not apply if the withdrawal is due to reasons specified within the section. Disclosure of estimated cut-off yield to eligible participants on the EBP platform is at the issuer''s discretion unless the withdrawal is due to a higher cut-off yield than the estimated one.
11.1: Issuers are required to open an escrow bank account or have a joint account with a Registrar, providing details of such accounts for pay-ins and setting timelines.
The section outlines regulatory responsibilities for parties involved in an issuance process. The Issuer is required to open an escrow bank account, provide details of this account to successful bidders, and credit securities after receiving confirmation of funds from the RTA. The RTA undertakes reconciliation between accounts, releases funds after passing instructions to depositories, and informs the Electronic Bidding Platform (EBP) about issue status.
The EBP is responsible for providing an online bidding platform, ensuring data safety and availability, and updating certain details on its website by specific deadlines. After May 16, 2025, the EBP must disclose estimated cut-off yields if an issuer withdraws due to a higher yield than the estimate. The section also emphasizes standardizing operational procedures across all EBP platforms and ensuring adherence to cooling-off periods and debarment periods.
The document outlines several regulations for Electronic Book Providers (EBPs) and electronic bid processes. Starting May 16, 2025, various clauses will be updated or revised, including provisions for investor grievances, bid limits, and allotment methods. EBPs must comply with these conditions, upgrade their systems, and communicate the changes to issuers, arrangers, and investors. Additionally, EBP platforms will be subjected to an annual audit by CISA. Chapter VII standardizes timelines for listing securities issued through private placement, requiring in-principle approval from stock exchanges before issuance and listing.
The section outlines a regulatory update affecting issuers looking to issue and list non-convertible securities or municipal debt securities. To receive in-principle approval from a stock exchange, they must submit a complete application following Regulation 6 of the NCS Regulations or Regulation 4A of the ILDM Regulations. The timelines for issuance and listing on private placement basis are provided in Table 1, with variations between EBP and Non-EBP activities. Key dates include receipt of in-principle approval before the issue open date (applicable from 3 months after May 16, 2025), bidding announcement on or before T-1, ISIN allocation prior to pay-in, finalization of allotments by the issuer, and listing application within timelines specified by the stock exchange(s).
Regulatory action: Starting from April 1, 2023, for private placement of debt securities, no more than 14 ISINs maturing in a financial year will be allowed per issuer. Additional six ISINs can be used for capital gains tax debt securities by authorized issuers under section 54EC of the Income Tax Act, 1961.
Affected parties: Issuers of debt securities and authorized issuers under section 54EC of the Income Tax Act, 1961.
Beginning April 1, 2023, for debt securities issued, the maximum number of ISINs maturing in a financial year is fourteen, with nine allocated for plain vanilla and structured debt securities, and five for market-linked debt securities. Issuers of only structured/market-linked debt securities are allowed up to nine ISINs per year. Previously, the limit was seventeen ISINs for all types of debt securities up to March 31, 2023. The limits may be reviewed periodically for potential reduction in fragmentation. Certain debt securities were previously given dispensations from ISIN restrictions until June 30, 2020, and the conversion of partly paid debt securities does not count as an additional ISIN.
- Regulatory amendment allows conversion of partly paid debt securities to fully paid without counting as an additional ISIN, applicable from December 31, 2023 (Section 4A).
- The limit of ISINs maturing in a financial year for unlisted ISINs converted to listed is excluded from the prescribed limit, also effective from December 31, 2023 (Section 4A).
- In case of exercised call/put options on debt securities, issuers may issue additional debt securities for the remaining period. The aggregate count of outstanding ISINs in the financial year should not exceed the prescribed limit (Section 5).
- Structured/ market linked debt securities maturity is reckoned based on original maturity date, not embedded options (Section 6).
- Debt obligations arising from ISIN restrictions can be honored through staggered or bullet maturity repayments (Section 7.1), different payment options for various investor categories (Section 7.2), and modifications in terms/structure of issue follow the procedure in Regulation 59 of SEBI LODR Regulations, 2015 (Section 9.3).
- Issuers must disclose the basis of payment when announcing multiple record dates (Section 9.4).
- An amendment to the AoA/charter/ constitution of the issuer is required for compliance with regulation 49(a) of SEBI NCS Regulations, 2021 to carry out consolidation and re-issuance of debt securities (Section 8).
- Issuers must submit a statement within fifteen days from the end of every half year to the stock exchange and depository containing data in a prescribed format (Section 9. Reporting and Monitoring, issuers).
The section outlines reporting and monitoring requirements for issuers of debt securities. Starting April 15 and October 15, issuers must submit a half-yearly report to stock exchanges and depositories in the specified format. Any modifications in issue terms should be communicated immediately. Stock exchanges and depositories are required to upload these reports on their respective websites and databases within prescribed timelines. Issuers who violate these provisions may be reported to SEBI within 30 working days of the half-year end. The document also includes an illustration for easier understanding, detailing the applicability of the circular for various issuance scenarios post April 1, 2023. Lastly, it mentions that stock exchanges and depositories must communicate the implementation status of these provisions to SEBI. This chapter is related to Green Debt Securities under Regulation 12A of SEBI NCS Regulations, 2021.
The regulatory update pertains to the implementation of Green Debt Securities under Chapter IX. This affects entities issuing or planning to issue green debt securities, as they must now disclose environmental sustainability objectives, project eligibility processes, and proposed deployment of proceeds in their offer documents. Additionally, within two years, they are required to appoint an independent third-party reviewer/certifier on a ''comply or explain'' basis. Post-issuance, listed green debt securities issuers must provide annual reports detailing the utilization of proceeds and unutilized amounts, with external auditor verification.
This regulatory update concerns Green Debt Securities. From April 1, 2023, issuers must disclose detailed information about unutilized proceeds and project allocations in their Annual Reports. Additionally, they are required to provide quantitative and qualitative environmental impact measures, methods used for preparation, mitigation plan deployment details, and environmental impact reporting on a project-by-project basis. A third-party reviewer/certifier is to be appointed for post-issue management and verification of internal tracking and impact reporting, with compliance expected within two years from the issue date. The issuer must also maintain a decision-making process for determining project eligibility, ensure compliance with greenwashing avoidance guidelines, and utilize proceeds only for stated purposes as disclosed in offer documents. Partial or non-compliance reasons and steps to achieve full compliance should be explained in the annual report. The BRSR essential indicators format is provided in Annexure IX-A.
The document outlines regulatory requirements for energy, water, and emissions reporting. Key actions include disclosing energy consumption (electricity, fuel, other sources), energy intensity, water withdrawal sources, and volumes. It also requires details on air emissions (NOx, SOx, PM, POP, VOC, HAP) and greenhouse gas emissions (Scope 1 and 2). Additionally, the document requests information about waste management, specifically total waste generated and recovered through recycling or reusing operations for various categories. An external agency assessment is optional for some sections. There is no mention of specific deadlines.
Regulatory body is seeking details on various waste categories (A-H) generated, their recovery methods, and disposal. It also requires information on waste management practices, usage of hazardous chemicals, and compliance with environmental laws in India, such as the Water Prevention and Control of Pollution Act, Air Prevention and Control of Pollution Act, and Environment Protection Act. Entities are asked to specify if any independent assessment has been conducted by an external agency, and to provide details on projects'' environmental impact assessments. Lastly, entities must provide a break-up of energy consumption from renewable and non-renewable sources and water discharge details with treatment levels. Compliance with these regulations is expected.
The document outlines regulatory requirements for reporting water discharge and withdrawal in areas of water stress, Scope 3 emissions, ecological impacts, and business continuity plans.
- Water Reporting: Facilities must provide details of water source, withdrawal, consumption, and discharge, including treatment levels (if any). Deadline not specified.
- Scope 3 Emissions: Entities must report total Scope 3 emissions and intensity in CO2 '
# Remove the now-obsolete row 3 (shrinks dimension to A1:L2)
$ws3.Rows(3).Delete()

# --- Sheet 2: was 'Press Release' -> 'Consulatation Paper' ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = 'Consulatation Paper'
$ws2.Range('A2').Value = 'SEBI'
$ws2.Range('B2').Value = 'Consulatation Paper'
$ws2.Range('C2').Value = 2025
$ws2.Range('D2').Value = 'October'
$ws2.Range('E2').NumberFormat = '@'
$ws2.Range('E2').Value = '2025-10-17'
$ws2.Range('E2').ClearFormats()
$ws2.Range('F2').Value = 'Consultation paper on proposed amendment to certain provisions of SEBI (LODR) Regulations, 2015 to facilitate transfer of securities transferred prior to April 1, 2019 and simplify the process of dematerialization of shares'
$ws2.Range('G2').Value = 'https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1760699641194.pdf'
$ws2.Range('H2').Value = '1760699641194.pdf'
$ws2.Range('I2').Value = '/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Consulatation Paper/2025/October/1760699641194.pdf'
$ws2.Range('J2').Value = 'No. 12/2019'
$ws2.Range('K2').Value = '"SEBI has proposed amendments to simplify the dematerialization process for securities transferred before April 1, 2019. A special window is open from July 7, 2025, to January 6, 2026, for re-lodging transfer deeds. The proposal also aims to eliminate the need for a Letter of Confirmation, reducing unnecessary steps and inconvenience for investors. Public comments are invited by November 7, 2025."'
$ws2.Range('L2').Value = 'SEBI has proposed an amendment to Regulation 40(1) regarding the transfer of securities, aimed at facilitating the transfer of securities transferred prior to April 1, 2019. This amendment is intended to simplify the process of dematerialization of securities. A special window for re-lodgement of transfer deeds has been opened from July 7, 2025, to January 6, 2026, to address cases where transfer deeds were previously rejected or returned due to deficiencies, and was extended due to investor requests and panel recommendations.
The document discusses a proposed amendment to Regulation 40(1) of LODR Regulations. This change aims to facilitate investors who executed transfer deeds before April 1, 2019, to get physical securities transferred and subsequently dematerialised. The proposal includes an exception for these investors, with a sunset clause to ensure maximum dematerialisation. A circular will be issued for the amendment''s implementation. Additionally, there''s a proposal to eliminate the need for a Letter of Confirmation (LOC) for investor service requests. This change would allow securities to be issued only in dematerialised form, with the RTA/listed company verifying and processing service requests and issuing LOCs within 30 days. The LOC is valid for 120 days from the date of issue. These changes are proposed to take effect after June 7, 2025, and will be in place until a specified sunset period.
The regulatory proposal aims to streamline the dematerialization process of securities by eliminating the issuance of Letter of Confirmation (LOC). Instead, depositories will develop a system for RTAs/listed entities to credit securities directly to investors'' demat accounts after due diligence. This change is intended to reduce unnecessary steps and inconvenience for investors, with a deadline for public comments on the proposal by November 7, 2025.
Proposal to amend certain provisions of SEBI (Listing Obligations and Disclosure Requirements) Regulations, 2015 is introduced. This affects the transfer of securities, including requirements for mandatory dematerialization, identity verification, signature verification, and advertisement for objections. Public comments are invited by 07-11-2025 through the provided link. In case of technical issues, emails can be sent to consultationMIRSD@sebi.gov.in. The proposal was issued on 17-10-2025.
Regulatory action: Amendment for Transfer-cum-Transmission cases where transfer deeds were executed before April 1, 2019 but not lodged, and the transferee has passed away.
Impacted parties: Listed entities, Registered Transfer Agents (RTA), claimants, and transferors.
Critical dates: Documents must be submitted by June 23, 2025; securities will be under lock-in for one year from the date of registration of transfer, and information related to transfers will be disclosed on company and stock exchange websites for six months from the date of transfer.
Regulatory document outlines an indemnity agreement between a company (M/s___________) and a buyer (I/We), effective under the Companies Act, 1956 and 2013. The agreement concerns securities purchased from sellers Mr. _______________ and ______________. To transfer securities into their name(s), the buyer has submitted the necessary documents with the Company/RTA. In exchange for this transfer approval, the buyer agrees to indemnify the company and its management, officers, staff, and the RTA from any harm, loss, damage, claims, demands, suits, actions, litigation, proceedings, costs, charges, or expenses that may arise as a result of the transfer. This indemnity bond is equivalent to the value of the shares at the prevailing market price with all accrued benefits. The agreement binds the buyer''s legal heirs, successors, administrators, and assigns. Signed by the Indemnifier in the presence of a Magistrate/Notary Public and witnesses on the specified date and place.'

